$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the intro blurb that currently lives in D2 - it is relocated to
# the new column K later on.
$introText = $ws.Range("D2").Value2

# Give E2 the same direct formatting (left/center/wrap) that D2 already has,
# before we overwrite the text, by copying D2's format onto E2.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# D2 and E2 are renamed to the new "mEFCT_PWRON" label (was the gun blurb /
# "mSPARKLE" respectively).
$ws.Range("D2").Value = "mEFCT_PWRON"
$ws.Range("E2").Value = "mEFCT_PWRON"

# The old D2 blurb moves into the new K2 cell, wrapped.
$ws.Range("K2").Value = $introText
$ws.Range("K2").WrapText = $true

# Give the new column K an explicit width (matches the author's ~30 char
# custom width).
$ws.Columns.Item(11).ColumnWidth = 29.8

# Leave the selection on D2:E2, matching where the edits were made.
$ws.Range("D2:E2").Select() | Out-Null
